$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "action_code_email_recipient" reference rows (40-42), continuing the
# pattern used by the other reference blocks in this lookup sheet.
$ws.Range("A40").Value = "action_code_email_recipient"
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "To Caller Person| Cc - Active personal in Email UPDT GENERAL GROUP  and  SLA group  | Bcc - Not Available"

$ws.Range("A41").Value = "action_code_email_recipient"
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = "To Selected personal in selected technician group| Cc - Other  active personal in the same selected technician group  | Bcc - Caller Person"

$ws.Range("A42").Value = "action_code_email_recipient"
$ws.Range("B42").Value = 3
$ws.Range("C42").Value = "To Selected personal in selected technician group| Cc - Other  active personal in the same selected technician group  | Bcc - Not Available"

# Widen the label/description columns to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 30.83
$ws.Columns.Item(3).ColumnWidth = 92

# Move the active selection the way the author left it after the edit.
$ws.Range("C43").Select() | Out-Null
